$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.7
$ws.Range("I2").Value = 2.25
$ws.Range("N2").Value = 5.6
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.57
$ws.Range("Q2").Value = 2.3
$ws.Range("S2").Value = 1.45
$ws.Range("T2").Value = 2.55
$ws.Range("U2").Value = 1.88
$ws.Range("V2").Value = 1.82
$ws.Range("W2").Value = 8.75
$ws.Range("X2").Value = 19.5
$ws.Range("Z2").Value = 60
$ws.Range("AC2").Value = 5.6
$ws.Range("AG2").Value = 700
$ws.Range("AK2").Value = 23
$ws.Range("AN2").Value = 5.5
$ws.Range("AP2").Value = 25
$ws.Range("AT2").Value = 2.55
$ws.Range("AU2").Value = 6.5
$ws.Range("AV2").Value = 55
$ws.Range("AX2").Value = 11.75
$ws.Range("AY2").Value = 19

$ws.Range("G5").Value = 1.98
$ws.Range("H5").Value = 4.1
$ws.Range("K5").Value = 2.62
$ws.Range("L5").Value = 3.2
$ws.Range("P5").Value = 6.5
$ws.Range("Q5").Value = 1.27
$ws.Range("R5").Value = 3.4
$ws.Range("T5").Value = 4.3
$ws.Range("U5").Value = 1.28
$ws.Range("V5").Value = 3.3
$ws.Range("W5").Value = 18
$ws.Range("X5").Value = 17
$ws.Range("Y5").Value = 9.75
$ws.Range("Z5").Value = 23
$ws.Range("AD5").Value = 10.5
$ws.Range("AF5").Value = 22
$ws.Range("AG5").Value = 80
$ws.Range("AH5").Value = 25
$ws.Range("AI5").Value = 29
$ws.Range("AJ5").Value = 12.5
$ws.Range("AK5").Value = 50
$ws.Range("AM5").Value = 18
$ws.Range("AQ5").Value = 28
$ws.Range("AR5").Value = 35
$ws.Range("AT5").Value = 4.3
$ws.Range("AV5").Value = 26
$ws.Range("AW5").Value = 6.2
$ws.Range("AX5").Value = 14.5
$ws.Range("AY5").Value = 14
$ws.Range("BB5").Value = 90

$ws.Range("AG6").Value = 800

$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 9
$ws.Range("Q7").Value = 1.8
$ws.Range("R7").Value = 2
